$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move from
# 45778 (2025-05-01) to 45779 (2025-05-02) for every data row (rows 2-43).
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45778) {
        $cell.Value = 45779
    }
}
